$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "61.036.73"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "3.382.60"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.77"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.93"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.68"
$ws.Range("E9").Value = "  +3.08%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "3.955.97"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.90"
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "3.373.25"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "61.058.80"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.89"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.22"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "76.11"
$ws.Range("E22").Value = "  +3.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.551"
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("E26").Value = "  +5.68%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.24"
$ws.Range("E28").Value = "  -1.57%  "
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.35"
$ws.Range("E32").Value = "  -3.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.32"
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.93"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "165.98"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").Value = "3.413.17"
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.29"
$ws.Range("E40").Value = "  -2.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.779"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.35"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("E44").Value = "  -2.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.12"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "2.452.35"
$ws.Range("E46").Value = "  -2.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.87"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.63"
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.13"
$ws.Range("E49").Value = "  +10.40%  "
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("E51").Value = "  -1.43%  "
